$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 51: date stored as text (matches existing rows' inlineStr style)
# and profit value, appended after the last existing data row (50).
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = "10/07/2025"
$ws.Range("B51").Value = 15208.78
